$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K9:K16").ClearContents()
$ws.Range("N9:N16").ClearContents()
